$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 97; everything from row 97 down (through 191)
# shifts down to 98..192, carrying its values/styles with it (old row 191's
# data lands on row 192, matching the diff's final hunk).
$ws.Rows("97:97").Insert()

# Populate the freshly-inserted, now-blank row 97 with the new weekly
# observation. Columns A,B,C,E,F,G,H,Q,R are constant across this whole
# market/category block, so copy them down from the row directly below
# (which holds the data that used to be row 97).
$ws.Cells.Item(97, 1).Value = 7
$ws.Cells.Item(97, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(97, 3).Value = "Ñuble"
$ws.Cells.Item(97, 4).Value = 44873
$ws.Cells.Item(97, 5).Value = 16
$ws.Cells.Item(97, 6).Value = 100112028
$ws.Cells.Item(97, 7).Value = "Sandia"
$ws.Cells.Item(97, 8).Value = "Sin especificar"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 300
$ws.Cells.Item(97, 11).Value = 800
$ws.Cells.Item(97, 12).Value = 900
$ws.Cells.Item(97, 13).Value = 850
$ws.Cells.Item(97, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(97, 15).Value = "Perú"
$ws.Cells.Item(97, 16).Value = 850
$ws.Cells.Item(97, 17).Value = 1
$ws.Cells.Item(97, 18).Value = "Hortaliza"
